# Refresh cryptocurrency Price (D) and Volume(1h) (E) columns with the
# latest scraped values (GitHub Actions cryptos-list update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "29.108.10"
$ws.Range('E2').Value = "  -0.96%  "
$ws.Range('D3').Value = "1.972.05"
$ws.Range('E3').Value = "  -0.80%  "
$ws.Range('D4').Value = "'1.013"
$ws.Range('E4').Value = "  +0.65%  "
$ws.Range('D5').Value = "'328.97"
$ws.Range('E5').Value = "  -0.28%  "
$ws.Range('D6').Value = "'1.010"
$ws.Range('E6').Value = "  +0.44%  "
$ws.Range('D7').Value = "'0.4954"
$ws.Range('E7').Value = "  -0.09%  "
$ws.Range('D8').Value = "'0.4205"
$ws.Range('E8').Value = "  +0.28%  "
$ws.Range('D9').Value = "'53.93"
$ws.Range('E9').Value = "  +4.08%  "
$ws.Range('D10').Value = "'0.09273"
$ws.Range('E10').Value = "  +4.86%  "
$ws.Range('D11').Value = "'1.099"
$ws.Range('E11').Value = "  -1.80%  "
$ws.Range('D12').Value = "'22.83"
$ws.Range('E12').Value = "  -1.83%  "
$ws.Range('D13').Value = "1.978.89"
$ws.Range('E13').Value = "  +0.70%  "
$ws.Range('D14').Value = "'7.891"
$ws.Range('E14').Value = "  -1.47%  "
$ws.Range('D15').Value = "'6.454"
$ws.Range('E15').Value = "  -0.56%  "
$ws.Range('E16').Value = "  +0.64%  "
$ws.Range('D17').Value = "'0.00001110"
$ws.Range('E17').Value = "  +0.60%  "
$ws.Range('D18').Value = "'91.89"
$ws.Range('E18').Value = "  -4.30%  "
$ws.Range('D19').Value = "'0.06723"
$ws.Range('E19').Value = "  +1.29%  "
$ws.Range('D20').Value = "'19.19"
$ws.Range('E20').Value = "  -2.40%  "
$ws.Range('E21').Value = "  +0.46%  "
$ws.Range('D22').Value = "'5.960"
$ws.Range('E22').Value = "  +0.24%  "
$ws.Range('D23').Value = "29.134.44"
$ws.Range('E23').Value = "  -0.90%  "
$ws.Range('D24').Value = "'11.97"
$ws.Range('E24').Value = "  +1.03%  "
$ws.Range('D25').Value = "'2.266"
$ws.Range('E25').Value = "  -0.67%  "
$ws.Range('D26').Value = "2.206.94"
$ws.Range('E26').Value = "  +0.20%  "
$ws.Range('D27').Value = "'20.76"
$ws.Range('E27').Value = "  +1.14%  "
$ws.Range('D28').Value = "'155.86"
$ws.Range('E28').Value = "  -0.93%  "
$ws.Range('D29').Value = "'6.231"
$ws.Range('E29').Value = "  -4.53%  "
$ws.Range('D30').Value = "'2.265"
$ws.Range('E30').Value = "  -3.05%  "
$ws.Range('D31').Value = "'127.17"
$ws.Range('E31').Value = "  -0.36%  "
$ws.Range('D32').Value = "'1.045"
$ws.Range('E32').Value = "  -0.42%  "
$ws.Range('D33').Value = "'0.09850"
$ws.Range('E33').Value = "  -0.65%  "
$ws.Range('D34').Value = "'1.500"
$ws.Range('E34').Value = "  -3.97%  "
$ws.Range('D35').Value = "'5.813"
$ws.Range('E35').Value = "  -0.38%  "
$ws.Range('D36').Value = "'3.730"
$ws.Range('E36').Value = "  -1.42%  "
$ws.Range('D37').Value = "'0.02424"
$ws.Range('E37').Value = "  -0.78%  "
$ws.Range('D38').Value = "'1.322"
$ws.Range('E38').Value = "  +3.19%  "
$ws.Range('D39').Value = "'0.06431"
$ws.Range('E39').Value = "  +1.43%  "
$ws.Range('D40').Value = "'9.032"
$ws.Range('E40').Value = "  -5.32%  "
$ws.Range('D41').Value = "'0.6482"
$ws.Range('E41').Value = "  -0.18%  "
$ws.Range('D42').Value = "'11.51"
$ws.Range('E42').Value = "  -1.94%  "
$ws.Range('D43').Value = "'0.2005"
$ws.Range('E43').Value = "  -2.89%  "
$ws.Range('E44').Value = "  +0.39%  "
$ws.Range('D45').Value = "'0.6213"
$ws.Range('E45').Value = "  -1.69%  "
$ws.Range('D46').Value = "'1.354"
$ws.Range('E46').Value = "  +7.30%  "
$ws.Range('D47').Value = "'13.31"
$ws.Range('E47').Value = "  -0.23%  "
$ws.Range('D48').Value = "'2.181"
$ws.Range('E48').Value = "  -1.11%  "
$ws.Range('D49').Value = "'3.485"
$ws.Range('E49').Value = "  -1.36%  "
$ws.Range('D50').Value = "'0.00000000328"
$ws.Range('E50').Value = "  +0.35%  "
$ws.Range('D51').Value = "'0.06973"
$ws.Range('E51').Value = "  -0.11%  "
